$wb = $excel.ActiveWorkbook

# --- Sheet "RunManager": zoom 80 -> 100 (zoomScale / zoomScaleNormal) ---
$wsRun = $wb.Worksheets.Item("RunManager")
$wsRun.Activate()
$excel.ActiveWindow.Zoom = 100

# --- Sheet "TestData": zoom 80 -> 100, add new row 8 (duplicate of row 7), move selection ---
$wsData = $wb.Worksheets.Item("TestData")
$wsData.Activate()
$excel.ActiveWindow.Zoom = 100

# Duplicate row 7 into row 8 via Copy so the shared-string cells (including the
# lone "'" values) keep their plain text style instead of picking up an
# auto "quote prefix" cell style that a direct Value/Formula assignment of "'" would trigger.
$wsData.Range("A7:G7").Copy($wsData.Range("A8:G8"))

# Move the active selection to C12, matching the saved view state.
$wsData.Range("C12").Select() | Out-Null
